# "2023 Day 18 done" -- refresh the running AoC leaderboard counts on the
# "2023" sheet (rows 2-18 get updated "My 1"/"My 2" counts as more people
# solve older puzzles), fill in the brand-new Day 18 row (row 19), and mark
# Day 17/Day 18 as completed on the "Overall" tracker sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "2023": updated leaderboard totals for days 1-17 (rows 2-18)
# ---------------------------------------------------------------------
$ws2023 = $wb.Worksheets.Item("2023")

$updates = @{
    2  = @(216683, 69033)
    3  = @(183246, 8513)
    4  = @(121230, 18235)
    5  = @(120677, 16233)
    6  = @(73869, 28633)
    7  = @(94739, 1598)
    8  = @(74619, 6602)
    9  = @(67172, 13171)
    10 = @(68619, 1027)
    11 = @(42776, 15347)
    12 = @(50683, 2068)
    13 = @(26250, 13262)
    14 = @(32659, 4568)
    15 = @(30670, 6480)
    16 = @(34024, 3587)
    17 = @(28286, 876)
    18 = @(16546, 1097)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws2023.Cells.Item($row, 2).Value = $vals[0]   # column B
    $ws2023.Cells.Item($row, 3).Value = $vals[1]   # column C
}

# New data for Day 18 (row 19) -- it was previously a blank template row.
$ws2023.Range("B19").Value = 13234
$ws2023.Range("C19").Value = 5408
$ws2023.Range("E19").Value = 12762
$ws2023.Range("F19").Value = 10338

# The row-19 formula cells were stored as blank string placeholders; give
# them back their real (shared) formulas so they evaluate now that the
# inputs are no longer blank.
$ws2023.Range("D19").Formula = '=IF(ISBLANK(B19),"",B19+C19)'
$ws2023.Range("G19").Formula = '=IF(D19="","",E19/D19)'
$ws2023.Range("H19").Formula = '=IF(ISBLANK(C19),"",F19/B19)'
$ws2023.Range("I19").Formula = '=IF(ISBLANK(E19),"",E19/$D$2)'
$ws2023.Range("J19").Formula = '=IF(ISBLANK(F19),"",F19/$B$2)'

# Selection left on B20 (next empty data row) after finishing Day 18 entry.
[void]$ws2023.Range("B20").Select()

# ---------------------------------------------------------------------
# Sheet "Overall": mark Day 17 / Day 18 (2023 row) as completed
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")

# Row 12 = year 2023. Columns BN:BQ = Day 17 (Test1,Test2,Live1,Live2),
# BR:BU = Day 18. All succeed ("s") except Day 17 Live-Part2 (BQ12) which
# is flagged as a performance issue ("p").
$wsOverall.Range("BN12").Value = "s"
$wsOverall.Range("BO12").Value = "s"
$wsOverall.Range("BP12").Value = "s"
$wsOverall.Range("BQ12").Value = "p"
$wsOverall.Range("BR12").Value = "s"
$wsOverall.Range("BS12").Value = "s"
$wsOverall.Range("BT12").Value = "s"
$wsOverall.Range("BU12").Value = "s"

# Selection moves to AY24, and "Overall" becomes the active tab/sheet.
[void]$wsOverall.Activate()
[void]$wsOverall.Range("AY24").Select()
